$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 24 (Test_Terrain): the terrain test is working again, so mark it Fixed
# and remove the old data-mismatch error note from column D.
$ws.Range("B24").Value = "Fixed"
$ws.Range("D24").ClearContents()

# Restore the previously selected cell
$ws.Range("C17").Select()

$wb.Save()
